$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 3) for the PF/1.0.7 entry, matching the existing
# columns: dev2 / sit2 / uat2 / prod
$ws.Range("A3").Value = "PF/1.0.7"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# The new row uses the sheet's default (unstyled) formatting rather than
# inheriting the style applied to rows 1-2.
$ws.Range("A3:D3").Style = "Normal"
